$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 2193
$ws.Range("F9").Value = 44
$ws.Range("F17").Value = 419
$ws.Range("F18").Value = 799
$ws.Range("F20").Value = 3028
$ws.Range("F22").Value = 122
$ws.Range("F23").Value = 3226
$ws.Range("F24").Value = 680
$ws.Range("F25").Value = 555
$ws.Range("F26").Value = 255
$ws.Range("F27").Value = 1006
$ws.Range("F28").Value = 743
$ws.Range("F30").Value = 789
$ws.Range("F31").Value = 770

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 150
$ws.Range("F20").Value = 205

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 420

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F13").Value = 420
$ws.Range("F14").Value = 2193
$ws.Range("F16").Value = 44
$ws.Range("F31").Value = 419
$ws.Range("F32").Value = 799
$ws.Range("F35").Value = 3028
$ws.Range("F36").Value = 122
$ws.Range("F37").Value = 3226
$ws.Range("F38").Value = 680
$ws.Range("F39").Value = 555
$ws.Range("F40").Value = 255
$ws.Range("F41").Value = 1006
$ws.Range("F44").Value = 205
$ws.Range("F45").Value = 146
$ws.Range("F47").Value = 743
$ws.Range("F49").Value = 789
$ws.Range("F50").Value = 770
